$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row 1 (shifts the former header row 1 -> row 2,
# and every data row down by one: old row N -> new row N+1).
$ws.Rows.Item(1).Insert()

# Row 1: a merged note/legend cell spanning the former header's columns.
$ws.Range("A1:D1").Merge()
$ws.Range("A1").Value = "Note: The date header (Row 2) supports: '2023 Annual', '2023 Q1', '2023-01'"
$ws.Range("A1").Font.Italic = $true
$ws.Range("A1").Font.Color = 255

# Row 2 (the former row 1 header, already carried its bold/centered style
# down when the new row was inserted) - update the year labels to include
# the period granularity ("Annual").
$ws.Range("B2").Value = "2024 Annual"
$ws.Range("C2").Value = "2023 Annual"
$ws.Range("D2").Value = "2022 Annual"
